$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Ghi chu" column (col E),
# shifting it (and the "Thuc hien" column after it) one to the right.
$ws.Columns.Item(5).Insert()

# Match the new column's width to its left neighbour (column D).
$ws.Columns.Item(5).ColumnWidth = 15

# --- Header cell (E2): "Tinh trang" ------------------------------------
# Copy column D's header formatting (bold, fill, centered, border) onto
# the new header cell, then set its text.
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E2").Value = "Tình trạng"

# --- Row 3 (Danh muc hang hoa): "Da hoan thanh" in green ---------------
$r3 = $ws.Range("E3")
$r3.Font.Color = 5287936
$r3.NumberFormat = "mm-dd-yy"
$r3.Value = "Đã hoàn thành"

# --- Row 4 (Nhap kho): "Dang thuc hien" ---------------------------------
$ws.Range("E4").Value = "Đang thực hiện"

# --- Row 5 (Tao user): "Dang thuc hien" ---------------------------------
$ws.Range("E5").Value = "Đang thực hiện"

# Rows 6-10 keep the inherited (empty) formatting from column D - no
# further changes required there.

# Update selection to match the post-edit cursor position.
$ws.Range("E5").Select()
